# Apply the "Add files via upload" edit: append 3 new spark-plug rows
# (ids 17-19) to the price sheet, turn the A1:E20 block into an AutoFilter
# range (with the usual hidden _FilterDatabase defined name Excel creates
# for it), and leave the selection on A19 like the author's workbook.
#
# Cells are written column-by-column (not row-by-row) so new shared
# strings land in the same order the source workbook has them in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New data rows (18, 19, 20 -> id 17, 18, 19)
# ---------------------------------------------------------------------

# A - id
$ws.Range("A18").Value = 17
$ws.Range("A19").Value = 18
$ws.Range("A20").Value = 19

# B - title
$ws.Range("B18").Value = "Свеча зажигания, Champion RCJ7Y"
$ws.Range("D18").Value = "RCJ7Y"
$ws.Range("B19").Value = "Свеча зажигания, Champion RC12YC"
$ws.Range("D19").Value = "RC12YC"
$ws.Range("B20").Value = "Свеча зажигания, Champion RJ19LM"
$ws.Range("D20").Value = "RJ19LM"

# E - manufacturer (reuses the existing "Champion" shared string)
$ws.Range("E18").Value = "Champion"
$ws.Range("E19").Value = "Champion"
$ws.Range("E20").Value = "Champion"

# G - description (reuses the existing "Данная деталь..." shared string)
$descr = "Данная деталь в наличии. Оплата товара за наличный расчет."
$ws.Range("G18").Value = $descr
$ws.Range("G19").Value = $descr
$ws.Range("G20").Value = $descr

# I - price
$ws.Range("I18").Value = 150
$ws.Range("I19").Value = 150
$ws.Range("I20").Value = 170

# K / L - delivery_days_from / delivery_days_to
$ws.Range("K18").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L18").Value = 1
$ws.Range("L19").Value = 1
$ws.Range("L20").Value = 1

# N - images
$ws.Range("N18").Value = "https://www.nixparts.com/assets/pictures/CHAMPION/RCJ7Y_01.JPG,https://www.nixparts.com/assets/pictures/CHAMPION/RCJ7Y_02.JPG"
$ws.Range("N19").Value = "https://www.nixparts.com/assets/pictures/CHAMPION/RC12YC_01.JPG,https://www.nixparts.com/assets/pictures/CHAMPION/RC12YC_02.JPG"
$ws.Range("N20").Value = "https://www.nixparts.com/assets/pictures/CHAMPION/RJ19LM_01.JPG,https://www.nixparts.com/assets/pictures/CHAMPION/RJ19LM_02.JPG"

# P - count
$ws.Range("P18").Value = 1
$ws.Range("P19").Value = 1
$ws.Range("P20").Value = 1

# ---------------------------------------------------------------------
# H and J columns (is_new / is_available) hold the literal text "True"
# in this workbook, not real booleans - assigning the string "True"
# directly gets auto-coerced to a boolean by the Value setter, so route
# it through a formula + paste-values round trip to keep it text, same
# as every other row in the sheet. (Done as two separate single-area
# ranges since .Formula on a multi-area union only touches the first
# area.)
# ---------------------------------------------------------------------
$rngH = $ws.Range("H18:H20")
$rngH.Formula = '=""&"True"'
$rngH.Copy()
$rngH.PasteSpecial(-4163)

$rngJ = $ws.Range("J18:J20")
$rngJ.Formula = '=""&"True"'
$rngJ.Copy()
$rngJ.PasteSpecial(-4163)

# ---------------------------------------------------------------------
# AutoFilter over the id/title/stores/part_number/manufacturer columns
# ---------------------------------------------------------------------
$ws.Range("A1:E20").AutoFilter()

# Excel auto-creates a hidden sheet-scoped _FilterDatabase name backing
# the AutoFilter - recreate it explicitly.
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "='15062018'!`$A`$1:`$E`$20")
$fd.Visible = $false

# ---------------------------------------------------------------------
# Final selection, matching the author's saved view
# ---------------------------------------------------------------------
$ws.Range("A19").Select()
